$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.940.12"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "2.300.23"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.15"
$ws.Range("E5").Value = "  +2.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.20"
$ws.Range("E6").Value = "  +3.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0968"
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.339"
$ws.Range("E11").Value = "  +4.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.91"
$ws.Range("E12").Value = "  +6.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.41"
$ws.Range("E13").Value = "  +7.91%  "
$ws.Range("D14").Value = "2.698.01"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "54.851.89"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000131"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("D17").Value = "2.286.86"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.42"
$ws.Range("E18").Value = "  +4.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.19"
$ws.Range("E19").Value = "  +3.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "308.18"
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.39"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "60.70"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.992"
$ws.Range("E24").Value = "  -2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.151"
$ws.Range("E25").Value = "  +2.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.48"
$ws.Range("E26").Value = "  +6.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "175.01"
$ws.Range("E27").Value = "  +5.81%  "
$ws.Range("D28").Value = "0.0₃0730"
$ws.Range("E28").Value = "  +7.95%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.09"
$ws.Range("E29").Value = "  +4.76%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.63"
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.11"
$ws.Range("E31").Value = "  +5.03%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.97"
$ws.Range("E33").Value = "  +2.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.992"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.947"
$ws.Range("E35").Value = "  +7.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.22"
$ws.Range("E36").Value = "  +4.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.81"
$ws.Range("E37").Value = "  +4.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.378"
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.44"
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.43"
$ws.Range("E40").Value = "  +2.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.00"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "126.66"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "253.12"
$ws.Range("E43").Value = "  +7.69%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0498"
$ws.Range("E44").Value = "  +3.46%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0901"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.554"
$ws.Range("E46").Value = "  +2.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.378"
$ws.Range("E47").Value = "  +2.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0209"
$ws.Range("E48").Value = "  +4.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.82"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.53"
$ws.Range("E50").Value = "  +3.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.58"
$ws.Range("E51").Value = "  +7.22%  "
